$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Week 11 semester long project entry (C13) to "Semester Long Project 3"
$ws.Range("C13").Value = "Semester Long Project 3"

# Update the active selection to C14
$ws.Range("C14").Select()
